# correction.json file grade format corrected
# Applies corrected Index/Name/sem1/sem2/sem3/CGPA values to the affected
# rows of the "CGPA" worksheet (rows shift as a consequence of the grade
# recalculation / re-sort that produced the corrected figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CGPA")

# Row => Index(B), Name(C), sem1(D), sem2(E), sem3(F), CGPA(G)
$rows = @(
    @{ Row = 24; B = 230258; C = "IMADUWAGE O.N.H.";          D = 3.9357; E = 3.964;  F = 3.9571; G = 3.955  },
    @{ Row = 25; B = 230508; C = "RAHUL B.";                  D = 4;      E = 4;      F = 3.8826; G = 3.9542 },
    @{ Row = 26; B = 230390; C = "MALDENIYA P.A.D.G.R.";      D = 4;      E = 4;      F = 3.8739; G = 3.9532 },
    @{ Row = 27; B = 230186; C = "FERNANDO W.H.D.";           D = 4;      E = 4;      F = 3.8696; G = 3.9516 },
    @{ Row = 28; B = 230159; C = "DISSANAYAKE G.R.G.K.";      D = 4;      E = 3.9571; F = 3.895;  G = 3.9468 },

    @{ Row = 32; B = 230197; C = "GARUSINGHE S.B.";           D = 4;      E = 3.928;  F = 3.9;    G = 3.9339 },
    @{ Row = 33; B = 230332; C = "KEERAWELLA K.P.C.P.";       D = 4;      E = 4;      F = 3.79;   G = 3.9323 },
    @{ Row = 34; B = 230486; C = "PERERA U.I.H.";             D = 4;      E = 4;      F = 3.8269; G = 3.9308 },
    @{ Row = 35; B = 230140; C = "DHARMAKEERTHI P.K.G.C.L.";  D = 3.9357; E = 3.964;  F = 3.8739; G = 3.9242 },
    @{ Row = 36; B = 230521; C = "RANASINGHE D.P.H.";         D = 4;      E = 4;      F = 3.7913; G = 3.9186 },
    @{ Row = 37; B = 230536; C = "RASANJANA W.P.G.R.A.";      D = 3.9571; E = 3.9591; F = 3.8478; G = 3.9153 },

    @{ Row = 94; B = 230013; C = "ABEYWARNA D.H.";            D = 3.85;   E = 3.648;  F = 3.6385; G = 3.6877 },
    @{ Row = 95; B = 230458; C = "PALIHENA H.H.";             D = 3.9571; E = 3.736;  F = 3.4652; G = 3.6855 },
    @{ Row = 96; B = 230735; C = "WITHANAGE G.D.N.";          D = 3.9357; E = 3.8182; F = 3.4;    G = 3.6831 },
    @{ Row = 97; B = 230248; C = "HIMASARA W.V.M.J.";         D = 3.9357; E = 3.684;  F = 3.5435; G = 3.6806 },
    @{ Row = 98; B = 230581; C = "SANTHOSH S.";               D = 3.7929; E = 3.684;  F = 3.5471; G = 3.6696 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
